# Append a new log entry (row 7) to the "Proximity" sheet, recording a
# "Bedroom Door" EXIT event, mirroring the existing ENTER/EXIT rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Proximity")

$row = 7

# Leading apostrophes force Excel to store these values as plain text
# instead of auto-converting "2026-02-01"/"15:06:26"/"15:00" into date
# or time values, matching how the rest of the sheet is stored.
$ws.Cells.Item($row, 1).Value = "'2026-02-01"
$ws.Cells.Item($row, 2).Value = "'15:06:26"
$ws.Cells.Item($row, 3).Value = "'15:00"
$ws.Cells.Item($row, 4).Value = "Bedroom Door"
$ws.Cells.Item($row, 5).Value = "EXIT"
$ws.Cells.Item($row, 6).Value = "User EXITED Bedroom"

# Resetting the style clears the "quote prefix" formatting Excel applies
# when a value is entered with a leading apostrophe, so the new cells end
# up with the same plain/default styling as the rest of the sheet.
$ws.Range($ws.Cells.Item($row, 1), $ws.Cells.Item($row, 6)).Style = "Normal"
